$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 699.1667
$ws.Range("I53").Value = 898.3333
$ws.Range("K53").Value = 898.3333
$ws.Range("M53").Value = -261.3333
$ws.Range("H125").Value = 1989.3572
$ws.Range("I125").Value = 797
$ws.Range("J125").Value = 3181.7144
$ws.Range("K125").Value = 7173
$ws.Range("L125").Value = 28635.4296
$ws.Range("M125").Value = -4713
$ws.Range("N125").Value = -33555.4296
$ws.Range("H138").Value = 16668423
$ws.Range("I138").Value = 27779230
$ws.Range("J138").Value = 2211.625
$ws.Range("K138").Value = 83337690
$ws.Range("L138").Value = 6634.875
$ws.Range("M138").Value = -83332550
$ws.Range("N138").Value = -16914.875

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3211.4443
$ws.Range("I2").Value = 2352.5
$ws.Range("J2").Value = 4929.3335
$ws.Range("K2").Value = 2352.5
$ws.Range("L2").Value = 4929.3335
$ws.Range("M2").Value = -2239.5
$ws.Range("N2").Value = -5155.3335
$ws.Range("H32").Value = 3239.1348
$ws.Range("I32").Value = 2866.5823
$ws.Range("J32").Value = 6182.3
$ws.Range("K32").Value = 2866.5823
$ws.Range("L32").Value = 6182.3
$ws.Range("M32").Value = -2579.5823
$ws.Range("N32").Value = -6756.3
$ws.Range("H45").Value = 6406.1816
$ws.Range("I45").Value = 8282.1
$ws.Range("J45").Value = 3520.1538
$ws.Range("K45").Value = 8282.1
$ws.Range("L45").Value = 3520.1538
$ws.Range("M45").Value = -7905.1
$ws.Range("N45").Value = -4274.1538
$ws.Range("H61").Value = 5189.911
$ws.Range("I61").Value = 4694.9
$ws.Range("K61").Value = 4694.9
$ws.Range("M61").Value = -4482.9
$ws.Range("H74").Value = 9295.179
$ws.Range("I74").Value = 9470.673000000001
$ws.Range("J74").Value = 8490.833000000001
$ws.Range("K74").Value = 9470.673000000001
$ws.Range("L74").Value = 8490.833000000001
$ws.Range("M74").Value = -8596.673000000001
$ws.Range("N74").Value = -10238.833
$ws.Range("H77").Value = 9295.179
$ws.Range("I77").Value = 9470.673000000001
$ws.Range("J77").Value = 8490.833000000001
$ws.Range("K77").Value = 47353.36500000001
$ws.Range("L77").Value = 42454.165
$ws.Range("M77").Value = -42985.36500000001
$ws.Range("N77").Value = -51190.165
$ws.Range("H110").Value = 4208.6665
$ws.Range("I110").Value = 3922.25
$ws.Range("J110").Value = 6500
$ws.Range("K110").Value = 3922.25
$ws.Range("L110").Value = 6500
$ws.Range("M110").Value = -1877.25
$ws.Range("N110").Value = -10590
$ws.Range("H116").Value = 3211.4443
$ws.Range("I116").Value = 2352.5
$ws.Range("J116").Value = 4929.3335
$ws.Range("K116").Value = 2352.5
$ws.Range("L116").Value = 4929.3335
$ws.Range("M116").Value = -58.5
$ws.Range("N116").Value = -9517.333500000001
$ws.Range("H122").Value = 9049.929
$ws.Range("H136").Value = 5189.911
$ws.Range("I136").Value = 4694.9
$ws.Range("K136").Value = 14084.7
$ws.Range("M136").Value = -11534.7

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3211.4443
$ws.Range("I3").Value = 2352.5
$ws.Range("J3").Value = 4929.3335
$ws.Range("K3").Value = 2352.5
$ws.Range("L3").Value = 4929.3335
$ws.Range("M3").Value = -2238.5
$ws.Range("N3").Value = -5157.3335
$ws.Range("H98").Value = 50000
$ws.Range("J98").Value = 50000
$ws.Range("L98").Value = 50000
$ws.Range("N98").Value = -55990
$ws.Range("H139").Value = 140936.44
$ws.Range("J139").Value = 146665.53
$ws.Range("L139").Value = 146665.53
$ws.Range("N139").Value = -156945.53

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 45.07143
$ws.Range("I7").Value = 53.555557
$ws.Range("J7").Value = 29.8
$ws.Range("K7").Value = 53.555557
$ws.Range("L7").Value = 29.8
$ws.Range("M7").Value = 59.444443
$ws.Range("N7").Value = -255.8
$ws.Range("H17").Value = 1850
$ws.Range("I17").Value = 1850
$ws.Range("K17").Value = 1850
$ws.Range("M17").Value = -1676
$ws.Range("H140").Value = 120218.375
$ws.Range("J140").Value = 120218.375
$ws.Range("L140").Value = 120218.375
$ws.Range("N140").Value = -130578.375

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 2201750.5
$ws.Range("J9").Value = 2917.3333
$ws.Range("L9").Value = 8751.999899999999
$ws.Range("N9").Value = -9199.999899999999
$ws.Range("H51").Value = 1058.4
$ws.Range("J51").Value = 2100
$ws.Range("L51").Value = 6300
$ws.Range("N51").Value = -7220
$ws.Range("H57").Value = 5443.6665
$ws.Range("I57").Value = 5000
$ws.Range("J57").Value = 5998.25
$ws.Range("K57").Value = 15000
$ws.Range("L57").Value = 17994.75
$ws.Range("M57").Value = -14441
$ws.Range("N57").Value = -19112.75
$ws.Range("H62").Value = 13166.667
$ws.Range("J62").Value = 15000
$ws.Range("L62").Value = 45000
$ws.Range("N62").Value = -46372
$ws.Range("H65").Value = 13166.667
$ws.Range("J65").Value = 15000
$ws.Range("L65").Value = 135000
$ws.Range("N65").Value = -141864
$ws.Range("H97").Value = 199.71428
$ws.Range("J97").Value = 229.4
$ws.Range("L97").Value = 688.2
$ws.Range("N97").Value = -1680.2
$ws.Range("H102").Value = 6997.6
$ws.Range("I102").Value = 4988
$ws.Range("J102").Value = 7500
$ws.Range("K102").Value = 14964
$ws.Range("L102").Value = 22500
$ws.Range("M102").Value = -12530
$ws.Range("N102").Value = -27368
$ws.Range("H121").Value = 1759.5
$ws.Range("I121").Value = 850
$ws.Range("J121").Value = 1842.1818
$ws.Range("K121").Value = 2550
$ws.Range("L121").Value = 5526.5454
$ws.Range("M121").Value = -1240
$ws.Range("N121").Value = -8146.5454
$ws.Range("H140").Value = 1078.4615
$ws.Range("I140").Value = 1078.4615
$ws.Range("K140").Value = 3235.3845
$ws.Range("M140").Value = 1944.6155

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 22442.25
$ws.Range("J38").Value = 24884.5
$ws.Range("L38").Value = 24884.5
$ws.Range("N38").Value = -25810.5
$ws.Range("H57").Value = 18333.334
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H102").Value = 6784.8125
$ws.Range("I102").Value = 9260.625
$ws.Range("K102").Value = 9260.625
$ws.Range("M102").Value = -7638.625
$ws.Range("H122").Value = 3405.1292
$ws.Range("I122").Value = 3090.6538
$ws.Range("K122").Value = 9271.9614
$ws.Range("M122").Value = -6821.9614

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 224666.2
$ws.Range("I5").Value = 224666.2
$ws.Range("K5").Value = 224666.2
$ws.Range("M5").Value = -224553.2
$ws.Range("H22").Value = 956.1
$ws.Range("I22").Value = 935.1667
$ws.Range("J22").Value = 987.5
$ws.Range("K22").Value = 935.1667
$ws.Range("L22").Value = 987.5
$ws.Range("M22").Value = -640.1667
$ws.Range("N22").Value = -1577.5
$ws.Range("H27").Value = 956.1
$ws.Range("I27").Value = 935.1667
$ws.Range("J27").Value = 987.5
$ws.Range("K27").Value = 935.1667
$ws.Range("L27").Value = 987.5
$ws.Range("M27").Value = -828.1667
$ws.Range("N27").Value = -1201.5
$ws.Range("H40").Value = 2851.6191
$ws.Range("I40").Value = 2675.8823
$ws.Range("J40").Value = 3598.5
$ws.Range("K40").Value = 2675.8823
$ws.Range("L40").Value = 3598.5
$ws.Range("M40").Value = -2539.8823
$ws.Range("N40").Value = -3870.5
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H61").Value = 1675
$ws.Range("I61").Value = 1675
$ws.Range("K61").Value = 1675
$ws.Range("M61").Value = -1473
$ws.Range("H82").Value = 1950.44
$ws.Range("I82").Value = 1946.3673
$ws.Range("J82").Value = 2150
$ws.Range("K82").Value = 1946.3673
$ws.Range("L82").Value = 2150
$ws.Range("M82").Value = -1585.3673
$ws.Range("N82").Value = -2872
$ws.Range("H85").Value = 1950.44
$ws.Range("I85").Value = 1946.3673
$ws.Range("J85").Value = 2150
$ws.Range("K85").Value = 1946.3673
$ws.Range("L85").Value = 2150
$ws.Range("M85").Value = -698.3672999999999
$ws.Range("N85").Value = -4646
$ws.Range("H113").Value = 1675
$ws.Range("I113").Value = 1675
$ws.Range("K113").Value = 1675
$ws.Range("M113").Value = 495
$ws.Range("H122").Value = 4051.5789
$ws.Range("I122").Value = 4256.6924
$ws.Range("K122").Value = 12770.0772
$ws.Range("M122").Value = -10320.0772
$ws.Range("H139").Value = 122828.8
$ws.Range("J139").Value = 122828.8
$ws.Range("L139").Value = 122828.8
$ws.Range("N139").Value = -133108.8

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 100000
$ws.Range("I21").Value = 100000
$ws.Range("K21").Value = 100000
$ws.Range("M21").Value = -99765
$ws.Range("H29").Value = 28799.6
$ws.Range("I29").Value = 28799.6
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 28799.6
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -28509.6
$ws.Range("N29").ClearContents()
$ws.Range("H35").Value = 100000
$ws.Range("I35").Value = 100000
$ws.Range("K35").Value = 100000
$ws.Range("M35").Value = -99710
$ws.Range("H37").Value = 49999
$ws.Range("I37").Value = 49999
$ws.Range("K37").Value = 49999
$ws.Range("M37").Value = -49796
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H126").Value = 8214.272000000001
$ws.Range("I126").Value = 5150.222
$ws.Range("K126").Value = 15450.666
$ws.Range("M126").Value = -12980.666
